$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1
$ws.Range("B7").Formula = "=fluid_mass*c_water*(set_temp-ambient_temp)"
$ws.Range("B7").Select()
